$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The backend service model & controller were implemented, so the "TODO"
# markers in column I for the /service endpoint rows (45-54) are no longer
# needed and are cleared out.
$ws.Range("I45:I54").ClearContents()

# Reset the view: scroll back to the top-left corner and select cell A1
# (previously the sheet was scrolled to A24 with A39:XFD40 selected).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select() | Out-Null
